$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow programmatic updates to the
# model-holdings figures and the "as of" date disclaimer below.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A59).
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-26 for illustrative purposes only and are subject to change."

# Refresh the daily return (D) and risk/alpha (E) figures for each holding row.
$ws.Range("D2").Value = 0.02360427416091257
$ws.Range("E2").Value = -0.008496923527688272
$ws.Range("D3").Value = 0.01790134873662645
$ws.Range("E3").Value = 0.003973773097556199
$ws.Range("D4").Value = 0.01834357690264815
$ws.Range("E4").Value = 0.01421923474663922
$ws.Range("D5").Value = 0.02020219983812631
$ws.Range("E5").Value = 0.01788963007883582
$ws.Range("D6").Value = 0.0193110251157144
$ws.Range("E6").Value = 0.0147347740667978
$ws.Range("D7").Value = 0.02696267894630347
$ws.Range("E7").Value = 0.002828854314002704
$ws.Range("D8").Value = 0.01928178035763968
$ws.Range("E8").Value = 0.004201680672268893
$ws.Range("D9").Value = 0.0191608493310064
$ws.Range("E9").Value = 0.01848032340565964
$ws.Range("D10").Value = 0.01861389331512255
$ws.Range("E10").Value = 0.01936305732484067
$ws.Range("D11").Value = 0.01928672035055771
$ws.Range("E11").Value = 0.02243737513447064
$ws.Range("D12").Value = 0.01901581113893315
$ws.Range("E12").Value = 0.01608579088471851
$ws.Range("D13").Value = 0.02054483774692061
$ws.Range("E13").Value = 0.01275343361674297
$ws.Range("D14").Value = 0.01946989528795812
$ws.Range("E14").Value = 0.003836317135549994
$ws.Range("D15").Value = 0.01796260464880998
$ws.Range("E15").Value = 0.008228460793804571
$ws.Range("D16").Value = 0.01805350051850166
$ws.Range("E16").Value = 0.01451337507114392
$ws.Range("D17").Value = 0.01547600981359233
$ws.Range("E17").Value = 0.01825842696629199
$ws.Range("D18").Value = 0.01553410413030832
$ws.Range("E18").Value = 0.0275014628437682
$ws.Range("D19").Value = 0.0166671409059868
$ws.Range("E19").Value = 0.01356285863328099
$ws.Range("D20").Value = 0.01967697979108177
$ws.Range("E20").Value = 0.01164892548704555
$ws.Range("D21").Value = 0.02022037901206465
$ws.Range("E21").Value = 0.02269129287598926
$ws.Range("D22").Value = 0.02073433587525609
$ws.Range("E22").Value = 0.02714164546225617
$ws.Range("D23").Value = 0.01855164940435542
$ws.Range("E23").Value = -0.004313788145071218
$ws.Range("D24").Value = 0.02099496990161115
$ws.Range("E24").Value = 0.02268235294117638
$ws.Range("D25").Value = 0.02129650706932747
$ws.Range("E25").Value = 0.02805819477434679
$ws.Range("D26").Value = 0.0204254875180211
$ws.Range("E26").Value = 0.03325980961225894
$ws.Range("D27").Value = 0.01923356602675975
$ws.Range("E27").Value = 0.02173913043478271
$ws.Range("D28").Value = 0.02745904943470673
$ws.Range("E28").Value = 0.0353835193540728
$ws.Range("D29").Value = 0.01896957280522043
$ws.Range("E29").Value = 0.03374999999999995
$ws.Range("D30").Value = 0.01274379853049043
$ws.Range("E30").Value = 0.01108647450110856
$ws.Range("D31").Value = 0.009371166565495612
$ws.Range("E31").Value = 0.01918819188191878
$ws.Range("D32").Value = 0.01642567405215368
$ws.Range("E32").Value = 0.02676659528907921
$ws.Range("D33").Value = 0.01928276835622329
$ws.Range("E33").Value = 0.001188707280832091
$ws.Range("D34").Value = 0.01834278650378127
$ws.Range("E34").Value = 0.01189296333002976
$ws.Range("D35").Value = 0.01892135847434049
$ws.Range("E35").Value = 0.07882534775888717
$ws.Range("D36").Value = 0.01731092078306397
$ws.Range("E36").Value = 0.02972399150743099
$ws.Range("D37").Value = 0.01953036080127476
$ws.Range("E37").Value = -0.002063983488132082
$ws.Range("D38").Value = 0.01956237195538357
$ws.Range("E38").Value = 0.002222222222222126
$ws.Range("D39").Value = 0.02555833775956699
$ws.Range("E39").Value = 0.006556160316674831
$ws.Range("D40").Value = 0.01606880896375547
$ws.Range("E40").Value = -0.01214953271028019
$ws.Range("D41").Value = 0.02165060576169158
$ws.Range("E41").Value = 0.01139018691588789
$ws.Range("D42").Value = 0.01938018501656676
$ws.Range("E42").Value = 0.01569159240604412
$ws.Range("D43").Value = 0.02006229923868781
$ws.Range("E43").Value = 0.008391608391608463
$ws.Range("D44").Value = 0.01771106020942409
$ws.Range("E44").Value = 0.003079291762894565
$ws.Range("D45").Value = 0.02079737018489011
$ws.Range("E45").Value = 0.026603325415677
$ws.Range("D46").Value = 0.01956671914915143
$ws.Range("E46").Value = 0.00977560542101763
$ws.Range("D47").Value = 0.0183412057060475
$ws.Range("E47").Value = 0.009276018099547478
$ws.Range("D48").Value = 0.01595815312239169
$ws.Range("E48").Value = 0.0222882615156017
$ws.Range("D49").Value = 0.01794284467713787
$ws.Range("E49").Value = -0.0121580547112462
$ws.Range("D50").Value = 0.01721488732073754
$ws.Range("E50").Value = 0.0103305785123966
$ws.Range("D51").Value = 0.01640571648076486
$ws.Range("E51").Value = 0.03456790123456788
$ws.Range("D52").Value = 0.01847952550775223
$ws.Range("E52").Value = 0.01497005988023958
$ws.Range("D53").Value = 0.01554635531274503
$ws.Range("E53").Value = 0.01771823681936047
$ws.Range("D54").Value = 0.007691568973366721
$ws.Range("E54").Value = 0.01579961464354529
$ws.Range("D55").Value = 0.007178402509042164
$ws.Range("E55").Value = -0.003523452983924202
$ws.Range("D56").Value = 1
$ws.Range("E56").Value = 0.01522644137137341
